$d = $word.ActiveDocument

# --- Part 1: update the date field result text from "19. September 2013"
# (split across two runs with a "_GoBack" bookmark in between) to a single
# run containing "26. September 2013" (the bookmark at this location is
# removed, since it is re-created further below at its new location).
$d.Content.Find.Execute("19. September 2013", $true, $false, $false, $false, $false, $true, 1, $false, "26. September 2013", 2) | Out-Null

# --- Part 2: remove " by two weeks" from the sentence about the deadline
# extension, and move the "_GoBack" bookmark to the new split point between
# "...the deadline" and ". Enclosed you find...".
$rng = $d.Content
$rng.Find.Execute(" of the Journal of Chromatography B and for your very kind extension of the deadline by two weeks. Enclosed you f", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $rng.Start
$end = $rng.End

$splitPos = $start + 84

# Insert the bookmark at the new split point (right after "...the deadline").
$bmRng = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# Insert a temporary bookmark right after the found text, so that deleting
# " by two weeks" below does not let the trailing ". Enclosed you f" run
# merge back together with the following "ind our manuscript..." run.
$tmpRng = $d.Range($end, $end)
$d.Bookmarks.Add("_TempBreak", $tmpRng)

# Delete " by two weeks".
$delRng = $d.Range($splitPos, $splitPos + 13)
$delRng.Text = ""

# Remove the temporary helper bookmark again.
$d.Bookmarks("_TempBreak").Delete()
